$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This sheet is a weekly price log for "Locoto" at Vega Modelo de
# Temuco. The update:
#   1. Removes the old row that held D=44669 (2022-05-04) data.
#   2. Inserts a brand-new row (new row 25) with a fresh weekly
#      observation (2022-09-14).
#   3. Inserts another brand-new row (new row 38) with a fresh
#      weekly observation (2022-09-12).
# Net effect: the sheet grows from 43 to 44 data rows.
# ------------------------------------------------------------------

# Step 1: delete the row holding the obsolete observation (original row 39,
# D=44669 / J=60 / K=L=6250 / M=6250).
$ws.Rows.Item(39).Delete()

# Step 2: insert a new blank row at row 25 (shifts old rows 25.. down by one).
$ws.Rows.Item(25).Insert()

# Step 3: insert another new blank row at row 38 (after the delete+first
# insert, the old row 37 now sits at row 38 - push it, and everything below,
# down by one more).
$ws.Rows.Item(38).Insert()

# ------------------------------------------------------------------
# Fill the two brand-new rows with their full records.
# ------------------------------------------------------------------

function Set-LocotoRow($r, $fecha, $volumen, $precio) {
    $ws.Cells.Item($r, 1).Value = 10
    $ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value = "La Araucanía"
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = 9
    $ws.Cells.Item($r, 6).Value = 100112042
    $ws.Cells.Item($r, 7).Value = "Locoto"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 9).Value = "Primera"
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $precio
    $ws.Cells.Item($r, 12).Value = $precio
    $ws.Cells.Item($r, 13).Value = $precio
    $ws.Cells.Item($r, 14).Value = "$/kilo"
    $ws.Cells.Item($r, 15).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($r, 16).Value = $precio
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

# New row 25: Fecha 2022-09-14 (serial 44818), Volumen 35, Precio 2700
Set-LocotoRow 25 44818 35 2700

# New row 38: Fecha 2022-09-12 (serial 44816), Volumen 90, Precio 2700
Set-LocotoRow 38 44816 90 2700

# Make sure the date cells keep the workbook's date number format (style
# index 2, used by every other cell in column D).
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(24, 4).NumberFormat
$ws.Cells.Item(38, 4).NumberFormat = $ws.Cells.Item(37, 4).NumberFormat
